$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.00021951
$ws.Range("F2").Value = 0.01400598
$ws.Range("G2").Value = 0.00039225899699999996

$ws.Range("E3").Value = 0.00196704
$ws.Range("F3").Value = 0.01721106
$ws.Range("G3").Value = 0.0024532412727272727

$ws.Range("E4").Value = 0.00419832
$ws.Range("F4").Value = 0.16120332
$ws.Range("G4").Value = 0.005188931775700934
